$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value. Cells whose new text would otherwise be
# auto-coerced to a number/date by Excel (plain decimals like "0.4902") are flagged
# so we force the Text number format first, preserving the original inline-string
# semantics from the source data (prices/links/volume deltas are all text).
$changes = @(
    @{ Cell = 'D2'; Value = '26.412.16'; NumericRisk = $false }
    @{ Cell = 'E2'; Value = '  -0.30%  '; NumericRisk = $false }
    @{ Cell = 'D3'; Value = '1.725.09'; NumericRisk = $false }
    @{ Cell = 'E3'; Value = '  -0.07%  '; NumericRisk = $false }
    @{ Cell = 'E4'; Value = '  +0.08%  '; NumericRisk = $false }
    @{ Cell = 'E5'; Value = '  -0.37%  '; NumericRisk = $false }
    @{ Cell = 'E6'; Value = '  +0.06%  '; NumericRisk = $false }
    @{ Cell = 'D7'; Value = '0.4902'; NumericRisk = $true }
    @{ Cell = 'E7'; Value = '  +1.89%  '; NumericRisk = $false }
    @{ Cell = 'D8'; Value = '0.2607'; NumericRisk = $true }
    @{ Cell = 'E8'; Value = '  -2.35%  '; NumericRisk = $false }
    @{ Cell = 'D9'; Value = '0.06201'; NumericRisk = $true }
    @{ Cell = 'E9'; Value = '  +0.20%  '; NumericRisk = $false }
    @{ Cell = 'D10'; Value = '1.723.63'; NumericRisk = $false }
    @{ Cell = 'E10'; Value = '  -0.16%  '; NumericRisk = $false }
    @{ Cell = 'D11'; Value = '0.07015'; NumericRisk = $true }
    @{ Cell = 'E11'; Value = '  -2.46%  '; NumericRisk = $false }
    @{ Cell = 'E12'; Value = '  -0.34%  '; NumericRisk = $false }
    @{ Cell = 'D13'; Value = '4.582'; NumericRisk = $true }
    @{ Cell = 'E13'; Value = '  +1.23%  '; NumericRisk = $false }
    @{ Cell = 'D14'; Value = '0.5998'; NumericRisk = $true }
    @{ Cell = 'E14'; Value = '  -1.90%  '; NumericRisk = $false }
    @{ Cell = 'D15'; Value = '77.29'; NumericRisk = $true }
    @{ Cell = 'E15'; Value = '  +0.20%  '; NumericRisk = $false }
    @{ Cell = 'D17'; Value = '26.417.82'; NumericRisk = $false }
    @{ Cell = 'E17'; Value = '  -0.32%  '; NumericRisk = $false }
    @{ Cell = 'E18'; Value = '  +0.06%  '; NumericRisk = $false }
    @{ Cell = 'D19'; Value = '0.000007156'; NumericRisk = $true }
    @{ Cell = 'E19'; Value = '  +3.03%  '; NumericRisk = $false }
    @{ Cell = 'D20'; Value = '11.34'; NumericRisk = $true }
    @{ Cell = 'E20'; Value = '  -1.63%  '; NumericRisk = $false }
    @{ Cell = 'D21'; Value = '1.944.05'; NumericRisk = $false }
    @{ Cell = 'D22'; Value = '4.471'; NumericRisk = $true }
    @{ Cell = 'E22'; Value = '  -1.02%  '; NumericRisk = $false }
    @{ Cell = 'D23'; Value = '8.592'; NumericRisk = $true }
    @{ Cell = 'E23'; Value = '  -2.26%  '; NumericRisk = $false }
    @{ Cell = 'D24'; Value = '5.163'; NumericRisk = $true }
    @{ Cell = 'E24'; Value = '  -1.56%  '; NumericRisk = $false }
    @{ Cell = 'D25'; Value = '137.55'; NumericRisk = $true }
    @{ Cell = 'E25'; Value = '  +0.38%  '; NumericRisk = $false }
    @{ Cell = 'D26'; Value = '15.22'; NumericRisk = $true }
    @{ Cell = 'E26'; Value = '  -0.67%  '; NumericRisk = $false }
    @{ Cell = 'D27'; Value = '1.393'; NumericRisk = $true }
    @{ Cell = 'D28'; Value = '106.88'; NumericRisk = $true }
    @{ Cell = 'E28'; Value = '  -0.35%  '; NumericRisk = $false }
    @{ Cell = 'E29'; Value = '  -4.20%  '; NumericRisk = $false }
    @{ Cell = 'D30'; Value = '3.957'; NumericRisk = $true }
    @{ Cell = 'E30'; Value = '  -0.12%  '; NumericRisk = $false }
    @{ Cell = 'D31'; Value = '0.07947'; NumericRisk = $true }
    @{ Cell = 'E31'; Value = '  -0.80%  '; NumericRisk = $false }
    @{ Cell = 'D32'; Value = '3.681'; NumericRisk = $true }
    @{ Cell = 'E32'; Value = '  -0.23%  '; NumericRisk = $false }
    @{ Cell = 'D33'; Value = '0.04534'; NumericRisk = $true }
    @{ Cell = 'E33'; Value = '  +0.25%  '; NumericRisk = $false }
    @{ Cell = 'B34'; Value = 'Frax'; NumericRisk = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; NumericRisk = $false }
    @{ Cell = 'D34'; Value = '0.9995'; NumericRisk = $true }
    @{ Cell = 'E34'; Value = '  +0.03%  '; NumericRisk = $false }
    @{ Cell = 'B35'; Value = 'HuobiToken'; NumericRisk = $false }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; NumericRisk = $false }
    @{ Cell = 'D35'; Value = '2.605'; NumericRisk = $true }
    @{ Cell = 'E35'; Value = '  -0.26%  '; NumericRisk = $false }
    @{ Cell = 'B36'; Value = 'ARBITRUM'; NumericRisk = $false }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; NumericRisk = $false }
    @{ Cell = 'D36'; Value = '0.9963'; NumericRisk = $true }
    @{ Cell = 'E36'; Value = '  -0.06%  '; NumericRisk = $false }
    @{ Cell = 'B37'; Value = 'ImmutableX'; NumericRisk = $false }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; NumericRisk = $false }
    @{ Cell = 'D37'; Value = '0.6279'; NumericRisk = $true }
    @{ Cell = 'E37'; Value = '  +0.14%  '; NumericRisk = $false }
    @{ Cell = 'B38'; Value = 'TrustWalletToken'; NumericRisk = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; NumericRisk = $false }
    @{ Cell = 'D38'; Value = '0.9089'; NumericRisk = $true }
    @{ Cell = 'E38'; Value = '  -0.39%  '; NumericRisk = $false }
    @{ Cell = 'B39'; Value = 'RenderToken'; NumericRisk = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; NumericRisk = $false }
    @{ Cell = 'D39'; Value = '1.957'; NumericRisk = $true }
    @{ Cell = 'E39'; Value = '  -5.86%  '; NumericRisk = $false }
    @{ Cell = 'B40'; Value = 'MXToken'; NumericRisk = $false }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; NumericRisk = $false }
    @{ Cell = 'D40'; Value = '2.393'; NumericRisk = $true }
    @{ Cell = 'E40'; Value = '  +1.03%  '; NumericRisk = $false }
    @{ Cell = 'B41'; Value = 'PaxDollar'; NumericRisk = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; NumericRisk = $false }
    @{ Cell = 'D41'; Value = '1.001'; NumericRisk = $true }
    @{ Cell = 'E41'; Value = '  -0.03%  '; NumericRisk = $false }
    @{ Cell = 'B42'; Value = 'VeChain'; NumericRisk = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; NumericRisk = $false }
    @{ Cell = 'D42'; Value = '0.01484'; NumericRisk = $true }
    @{ Cell = 'E42'; Value = '  -1.59%  '; NumericRisk = $false }
    @{ Cell = 'B43'; Value = 'Quant'; NumericRisk = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; NumericRisk = $false }
    @{ Cell = 'D43'; Value = '99.64'; NumericRisk = $true }
    @{ Cell = 'E43'; Value = '  -3.59%  '; NumericRisk = $false }
    @{ Cell = 'B44'; Value = 'FraxShare'; NumericRisk = $false }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; NumericRisk = $false }
    @{ Cell = 'D44'; Value = '5.445'; NumericRisk = $true }
    @{ Cell = 'E44'; Value = '  -3.49%  '; NumericRisk = $false }
    @{ Cell = 'B45'; Value = 'TheSandbox'; NumericRisk = $false }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; NumericRisk = $false }
    @{ Cell = 'D45'; Value = '0.3845'; NumericRisk = $true }
    @{ Cell = 'E45'; Value = '  -0.50%  '; NumericRisk = $false }
    @{ Cell = 'B46'; Value = 'Aptos'; NumericRisk = $false }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; NumericRisk = $false }
    @{ Cell = 'D46'; Value = '6.717'; NumericRisk = $true }
    @{ Cell = 'E46'; Value = '  -3.73%  '; NumericRisk = $false }
    @{ Cell = 'B47'; Value = 'Algorand'; NumericRisk = $false }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; NumericRisk = $false }
    @{ Cell = 'D47'; Value = '0.1158'; NumericRisk = $true }
    @{ Cell = 'E47'; Value = '  -2.00%  '; NumericRisk = $false }
    @{ Cell = 'B48'; Value = 'Cronos'; NumericRisk = $false }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; NumericRisk = $false }
    @{ Cell = 'D48'; Value = '0.05367'; NumericRisk = $true }
    @{ Cell = 'E48'; Value = '  +0.24%  '; NumericRisk = $false }
    @{ Cell = 'B49'; Value = 'EnergySwap'; NumericRisk = $false }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; NumericRisk = $false }
    @{ Cell = 'D49'; Value = '7.760'; NumericRisk = $true }
    @{ Cell = 'E49'; Value = '  -0.90%  '; NumericRisk = $false }
    @{ Cell = 'B50'; Value = 'Elrond'; NumericRisk = $false }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'; NumericRisk = $false }
    @{ Cell = 'D50'; Value = '30.13'; NumericRisk = $true }
    @{ Cell = 'E50'; Value = '  -1.09%  '; NumericRisk = $false }
    @{ Cell = 'B51'; Value = 'NEARProtocol'; NumericRisk = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; NumericRisk = $false }
    @{ Cell = 'D51'; Value = '1.240'; NumericRisk = $true }
    @{ Cell = 'E51'; Value = '  -0.93%  '; NumericRisk = $false }
)

foreach ($change in $changes) {
    $range = $ws.Range($change.Cell)
    if ($change.NumericRisk) {
        $range.NumberFormat = '@'
    }
    $range.Value = $change.Value
}
